# Applies the "adding averages and more checks" update:
#  - Bold+white font for the dashboard titles and table headers (both sheets)
#  - Training Dashboard: PERIOD TO EXPIRE (-8) and LAST UPDATE (08->16-Sep-2025)
#    refreshed for every training row
#  - Training Dashboard row 21 (LOTO SOP) flips to NOT VALID / red highlight,
#    matching the look already used by rows 24-25
#  - Exam Dashboard: widen COMMENTS column and reword "OK" -> "date is valid"

$wb = $excel.ActiveWorkbook

$white = 16777215

$wsTrain = $wb.Worksheets.Item("Training Dashboard")
$wsExam  = $wb.Worksheets.Item("Exam Dashboard")

# --- Header / title font styling (bold, white text) -------------------------
$wsTrain.Range("A1").Font.Bold = $true
$wsTrain.Range("A1").Font.Size = 11
$wsTrain.Range("A1").Font.Color = $white

$wsTrain.Range("A2:K2").Font.Bold = $true
$wsTrain.Range("A2:K2").Font.Size = 11
$wsTrain.Range("A2:K2").Font.Color = $white

$wsExam.Range("A1:G1").Font.Bold = $true
$wsExam.Range("A1:G1").Font.Size = 11
$wsExam.Range("A1:G1").Font.Color = $white

$wsExam.Range("A2:G2").Font.Bold = $true
$wsExam.Range("A2:G2").Font.Size = 11
$wsExam.Range("A2:G2").Font.Color = $white

# --- Training Dashboard: refresh "period to expire" / "last update" ---------
foreach ($row in 3..27) {
    $cell = $wsTrain.Range("H$row")
    $cell.Value2 = $cell.Value2 - 8
    $wsTrain.Range("I$row").Value2 = "'16-Sep-2025"
}

# --- Row 21 (LOTO SOP) becomes invalid, mirror rows 24/25 formatting --------
$wsTrain.Range("A24:K24").Copy()
$wsTrain.Range("A21:K21").PasteSpecial(-4122)

$wsTrain.Range("H21").Value2 = 15
$wsTrain.Range("I21").Value2 = "'16-Sep-2025"
$wsTrain.Range("J21").Value2 = "NOT VALID"

# --- Exam Dashboard: widen comments column, reword remark -------------------
$wsExam.Columns.Item(5).ColumnWidth = 14.14

foreach ($row in 3..10) {
    $wsExam.Range("E$row").Value2 = "date is valid"
}
